# Worker Marg Income Tax Rate.xlsx
# Update the "WMITR" sheet's column header (B1) to clarify the metric is
# dimensionless, enable word-wrap on that header cell, and let the header
# row grow to fit the now-longer, wrapped text. Keep the "About" sheet the
# active/selected sheet, just as it was before the edit.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("WMITR")

# Update the header text in B1 of the WMITR sheet.
$wsData.Range("B1").Value = "Marginal Income Tax Rate (dimensionless)"

# Turn on wrap text for that header cell (it was right-aligned already).
$wsData.Range("B1").WrapText = $true

# The longer, wrapped header now needs a taller row.
$wsData.Rows.Item(1).RowHeight = 28.5

# Leave B1 selected on the WMITR sheet ...
$wsData.Range("B1").Select()

# ... but restore "About" as the active sheet/tab, matching the original
# workbook state (tabSelected stayed on "About").
$wsAbout.Activate()
